$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "SA4"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "September 30, 2024"

$ws.Range("F10").Select()
